# edit.ps1 - applies the "add basic search bar" change set to the
# Ecommerce Application Design document.
#
# Substantive (visible) changes made here:
#   1. Highlight (yellow) the "Menu/Categories (get)" line.
#   2. Highlight (yellow) the "Menu/Categories/category :id (get)" line.
#   3. Highlight (yellow) the "Add search bar functionality" line.
#   4. Insert a new paragraph after it:
#        "Add autocomplete functionality & cleanup categories"
#   5. After "Create payment, delivery", insert three new paragraphs:
#        "Bugs:"
#        "Fix page resize"
#        "Fix Footer on main page"
#   6. Remove the stray <w:lastRenderedPageBreak/> that sits on the
#      trailing tab-only paragraph near the end of the document (the
#      one in the "Users: name, email..." paragraph must stay put).
#
# (The rest of the upstream diff only wraps already-correct words in
# <w:proofErr> spell-check bookmarks produced automatically by Word's
# background spell checker -- it does not change the document's visible
# text/formatting, so there is nothing more to do for those hunks.)

$d = $word.ActiveDocument

# --- 1. Highlight "Menu/Categories (get)" -------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Menu/Categories (get)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 7
}

# --- 2. Highlight "Menu/Categories/category :id (get)" ------------------
$rng = $d.Content
$found = $rng.Find.Execute("Menu/Categories/category :id (get)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.HighlightColorIndex = 7
}

# --- 3 & 4. Insert the new "Add autocomplete ..." paragraph right after -
#            "Add search bar functionality", then highlight only the
#            original line (inserting first, highlighting second, keeps
#            the new paragraph's runs free of any inherited <w:rPr>).
$rng = $d.Content
$found = $rng.Find.Execute("Add search bar functionality", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.InsertParagraphAfter()
    $searchBarPara = $rng.Paragraphs(1)
    $autocompletePara = $searchBarPara.Next()
    $autocompletePara.Range.Text = "Add autocomplete functionality & cleanup categories"
    $rng.HighlightColorIndex = 7
}

# --- 5. Insert "Bugs:" / "Fix page resize" / "Fix Footer on main page" --
#        after "Create payment, delivery".
$rng = $d.Content
$found = $rng.Find.Execute("Create payment, delivery", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.InsertParagraphAfter()
    $paymentPara = $rng.Paragraphs(1)

    $bugsPara = $paymentPara.Next()
    $bugsPara.Range.Text = "Bugs:"

    $bugsPara.Range.InsertParagraphAfter()
    $resizePara = $bugsPara.Next()
    $resizePara.Range.Text = "Fix page resize"

    $resizePara.Range.InsertParagraphAfter()
    $footerPara = $resizePara.Next()
    $footerPara.Range.Text = "Fix Footer on main page"
}

# --- 6. Remove the trailing stray lastRenderedPageBreak -----------------
# There are three consecutive tab-only paragraphs near the end of the
# document; the last of the three carries the spurious
# <w:lastRenderedPageBreak/> run child. Re-assigning its Range.Text
# rewrites the run cleanly without that element, while leaving the
# legitimate lastRenderedPageBreak (in the "Users: name, email..."
# paragraph earlier in the doc) untouched.
$tabOnlyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq [char]9) {
        $tabOnlyIndex = $i
    }
}
if ($tabOnlyIndex -ge 1) {
    $d.Paragraphs($tabOnlyIndex).Range.Text = [char]9
}

Write-Output "done"
